$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("swiss_matches")

# The swiss draw table now only has 3 match rows instead of 4 - drop the last one
$ws.Rows.Item(4).Delete()

# Set the updated match data
$ws.Range("A1").Value = "Cosmo"
$ws.Range("B1").Value = "Frizmi"
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 10

$ws.Range("A2").Value = "Vertigo"
$ws.Range("B2").Value = "Cosmo"
$ws.Range("C2").Value = 15
$ws.Range("D2").Value = 10

$ws.Range("A3").Value = "Vertigo"
$ws.Range("B3").Value = "Frizmi"
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = 15

$ws.Range("D3").Select()
